$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Find-ParagraphByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Trim() -eq $needle) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Edit 1 (Project Description paragraph): drop the "existing websites such as
# amazon.com and allrecipies.com..." sentence in favor of "a recipe API...".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "using existing websites such as amazon.com and allrecipies.com to use their platforms to achieve our website" + [char]0x2019 + "s purpose.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "using a recipe API to achieve our site" + [char]0x2019 + "s purpose.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 2 (Project Deliverables paragraph): merge "a recipe" + "_GoBack"
# bookmark + " database " into a single "a recipe database " run (the
# bookmark is relocated later, near the end of the document).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "a recipe database ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a recipe database ",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 3 (Project Milestones, first bullet): replace the web-scraping
# milestone description with the new API milestone.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Web Scraping Implementation- Develop the software to pull ingredients from allrecipies.com and search them on amazon.com.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "API implementation and testing.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 4 (Production Fees and Total Cost section): the page-break marker
# that used to render in front of "10 weeks total" now renders in front of
# "Software implementation - 50%" (the section grew by a paragraph), and
# three new paragraphs are appended at the end: a blank line, a Labor cost
# line, and a Web Server cost line (the relocated "_GoBack" bookmark sits
# between "Web " and "Server - $100").
# ---------------------------------------------------------------------------
$weeksPara = Find-ParagraphByText $d "10 weeks total"
$weeksPara.Range.InsertXML("<w:p xmlns:w='$wNs'><w:pPr><w:ind w:left='360'/></w:pPr><w:r><w:t>10 weeks total</w:t></w:r></w:p>") | Out-Null

$softwareNeedle = "Software implementation " + [char]0x2013 + " 50%"
$softwarePara = Find-ParagraphByText $d $softwareNeedle
$softwareXml = "<w:p xmlns:w='$wNs'><w:pPr><w:ind w:left='360'/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Software implementation " + [char]0x2013 + " 50%</w:t></w:r></w:p>"
$softwarePara.Range.InsertXML($softwareXml) | Out-Null

# Append: blank paragraph, Labor line, Web Server line - after the last
# existing paragraph ("Revising software- 20%").
$revisingPara = Find-ParagraphByText $d "Revising software- 20%"
$revisingPara.Range.InsertParagraphAfter() | Out-Null

$blankPara = $d.Paragraphs($d.Paragraphs.Count)
$blankPara.Range.InsertParagraphAfter() | Out-Null

$laborPara = $d.Paragraphs($d.Paragraphs.Count)
$laborPara.Range.Text = "Labor " + [char]0x2013 + " 60hrs x `$90 x 2 = `$10,800"
$laborPara.Range.InsertParagraphAfter() | Out-Null

$webPara = $d.Paragraphs($d.Paragraphs.Count)
$webPara.Range.Text = "Web Server - `$100"

# Re-insert the "_GoBack" bookmark between "Web " and "Server - $100".
$bmRange = $d.Range($webPara.Range.Start + 4, $webPara.Range.Start + 4)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
